$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values while preserving them as plain text (avoids Excel
# auto-converting numeric-looking strings like "1.12" into numbers).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '90.535.91'
Set-TextValue $ws.Range('E2') '  +0.49%  '
Set-TextValue $ws.Range('D3') '3.104.34'
Set-TextValue $ws.Range('E3') '  +0.95%  '
Set-TextValue $ws.Range('E4') '  -0.06%  '
Set-TextValue $ws.Range('D5') '239.60'
Set-TextValue $ws.Range('E5') '  +9.88%  '
Set-TextValue $ws.Range('D6') '623.77'
Set-TextValue $ws.Range('E6') '  +1.64%  '
Set-TextValue $ws.Range('D7') '1.12'
Set-TextValue $ws.Range('E7') '  +6.45%  '
Set-TextValue $ws.Range('D8') '0.370'
Set-TextValue $ws.Range('E8') '  +5.80%  '
Set-TextValue $ws.Range('E9') '  -0.01%  '
Set-TextValue $ws.Range('D10') '3.101.92'
Set-TextValue $ws.Range('E10') '  +1.03%  '
Set-TextValue $ws.Range('D11') '0.739'
Set-TextValue $ws.Range('E11') '  +3.22%  '
Set-TextValue $ws.Range('E12') '  +3.72%  '
Set-TextValue $ws.Range('D13') '0.0000249'
Set-TextValue $ws.Range('E13') '  +4.11%  '
Set-TextValue $ws.Range('D14') '35.12'
Set-TextValue $ws.Range('E14') '  +2.68%  '
Set-TextValue $ws.Range('D15') '5.47'
Set-TextValue $ws.Range('E15') '  -0.36%  '
Set-TextValue $ws.Range('D16') '90.436.63'
Set-TextValue $ws.Range('E16') '  +0.74%  '
Set-TextValue $ws.Range('D17') '3.693.07'
Set-TextValue $ws.Range('E17') '  +1.73%  '
Set-TextValue $ws.Range('D18') '3.085.43'
Set-TextValue $ws.Range('E18') '  +0.96%  '
Set-TextValue $ws.Range('E19') '  +3.71%  '
Set-TextValue $ws.Range('D20') '14.28'
Set-TextValue $ws.Range('E20') '  +0.88%  '
Set-TextValue $ws.Range('D21') '0.0000211'
Set-TextValue $ws.Range('E21') '  +4.03%  '
Set-TextValue $ws.Range('D22') '5.71'
Set-TextValue $ws.Range('E22') '  +5.83%  '
Set-TextValue $ws.Range('D23') '446.79'
Set-TextValue $ws.Range('E23') '  +0.43%  '
Set-TextValue $ws.Range('D24') '8.99'
Set-TextValue $ws.Range('E24') '  +1.69%  '
Set-TextValue $ws.Range('D25') '5.90'
Set-TextValue $ws.Range('E25') '  +2.77%  '
Set-TextValue $ws.Range('D26') '91.31'
Set-TextValue $ws.Range('E26') '  +0.89%  '
Set-TextValue $ws.Range('D27') '12.02'
Set-TextValue $ws.Range('E27') '  +1.47%  '
Set-TextValue $ws.Range('D28') '3.262.20'
Set-TextValue $ws.Range('E28') '  +1.23%  '
Set-TextValue $ws.Range('D30') '0.179'
Set-TextValue $ws.Range('E30') '  +13.07%  '
Set-TextValue $ws.Range('D31') '0.216'
Set-TextValue $ws.Range('E31') '  +7.57%  '
Set-TextValue $ws.Range('D32') '9.20'
Set-TextValue $ws.Range('E32') '  -0.81%  '
Set-TextValue $ws.Range('D33') '1.01'
Set-TextValue $ws.Range('E33') '  +0.69%  '
Set-TextValue $ws.Range('D34') '0.113'
Set-TextValue $ws.Range('E34') '  +32.20%  '
Set-TextValue $ws.Range('B35') 'MantraDAO'
Set-TextValue $ws.Range('C35') 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue $ws.Range('D35') '4.29'
Set-TextValue $ws.Range('E35') '  +43.89%  '
Set-TextValue $ws.Range('B36') 'EthereumClassic'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D36') '26.44'
Set-TextValue $ws.Range('E36') '  -5.12%  '
Set-TextValue $ws.Range('D37') '0.156'
Set-TextValue $ws.Range('E37') '  +4.48%  '
Set-TextValue $ws.Range('D38') '7.54'
Set-TextValue $ws.Range('E38') '  +11.07%  '
Set-TextValue $ws.Range('D39') '1.92'
Set-TextValue $ws.Range('E39') '  +2.06%  '
Set-TextValue $ws.Range('D40') '494.09'
Set-TextValue $ws.Range('E40') '  -0.36%  '
Set-TextValue $ws.Range('D41') '3.60'
Set-TextValue $ws.Range('E41') '  +5.32%  '
Set-TextValue $ws.Range('D42') '1.29'
Set-TextValue $ws.Range('E42') '  +1.33%  '
Set-TextValue $ws.Range('E43') '  -2.07%  '
Set-TextValue $ws.Range('D44') '22.11'
Set-TextValue $ws.Range('E44') '  -0.41%  '
Set-TextValue $ws.Range('E45') '  +0.00%  '
Set-TextValue $ws.Range('D46') '159.90'
Set-TextValue $ws.Range('E46') '  +8.06%  '
Set-TextValue $ws.Range('D47') '1.91'
Set-TextValue $ws.Range('E47') '  -1.39%  '
Set-TextValue $ws.Range('D48') '0.681'
Set-TextValue $ws.Range('E48') '  -0.64%  '
Set-TextValue $ws.Range('D49') '4.53'
Set-TextValue $ws.Range('E49') '  +0.04%  '
Set-TextValue $ws.Range('D50') '44.83'
Set-TextValue $ws.Range('E50') '  +0.72%  '
Set-TextValue $ws.Range('D51') '1.33'
Set-TextValue $ws.Range('E51') '  +1.04%  '
